# "adjust excel to storage"
# - Insert a new "Storage" worksheet between "Transmission" and "Demand"
#   with cost/capacity parameters for a Left/Hydrogen storage.
# - Fix Process sheet inst-cap for Left/Mage (E3) from 40 to 12.
# - Make "Process" the active sheet (was "Transmission").

$wb = $excel.ActiveWorkbook

# --- 1. Build the new "Storage" worksheet -----------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Storage"
$demandSheet = $wb.Worksheets.Item("Demand")
$newSheet.Move($demandSheet)

# Re-fetch by name: .Move() leaves the old variable pointing at a stale
# position, so grab a fresh reference before writing any data.
$storage = $wb.Worksheets.Item("Storage")

$headers = @("Site","Storage","inst-cap-c","cap-lo-c","cap-up-c","inst-cap-p","cap-lo-p","cap-up-p","eff-in","eff-out","inv-cost-p","inv-cost-c","fix-cost-p","fix-cost-c","var-cost-p","var-cost-c","depreciation","wacc","init")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $storage.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$row2 = @("Left","Hydrogen",0,0,10,0,0,"inf",1,1,1,1,1,1,0.1,0.1,50,1,0.5)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $storage.Cells.Item(2, $i + 1).Value = $row2[$i]
}

# --- 2. Process sheet: fix Left/Mage inst-cap (E3) 40 -> 12 -----------
$process = $wb.Worksheets.Item("Process")
$process.Range("E3").Value = 12

# --- 3. Active sheet moves from Transmission to Process ---------------
$process.Activate()
